$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add column D header "rename", matching style of B1/C1
$ws.Cells.Item(1, 4).Value = 'rename'
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Data rows 2-55: update/insert A (index), B (feat name), C (shap value), D (rename/latex label)
# Row 2: From_Same_Orbital
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 'From_Same_Orbital'
$ws.Cells.Item(2, 3).Value = 0.00001632916089462153
$ws.Cells.Item(2, 4).Value = '$\mathbf{b}$'
# Row 3: coulomb
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 'coulomb'
$ws.Cells.Item(3, 3).Value = 0.003634297862900921
$ws.Cells.Item(3, 4).Value = '$\langle pp \vert \vert qq \rangle$'
# Row 4: screen1_1
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 'screen1_1'
$ws.Cells.Item(4, 3).Value = 0.0001235668008816771
$ws.Cells.Item(4, 4).Value = '$(\langle pp \vert \vert rr \rangle)_{1}$'
# Row 5: screen1_2
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 'screen1_2'
$ws.Cells.Item(5, 3).Value = 0.0007883970718058197
$ws.Cells.Item(5, 4).Value = '$(\langle pp \vert \vert rr \rangle)_{2}$'
# Row 6: screen1_3
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 'screen1_3'
$ws.Cells.Item(6, 3).Value = 0.00006748331070874053
$ws.Cells.Item(6, 4).Value = '$(\langle pp \vert \vert rr \rangle)_{3}$'
# Row 7: screen1_4
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 'screen1_4'
$ws.Cells.Item(7, 3).Value = 0.0002619603107121422
$ws.Cells.Item(7, 4).Value = '$(\langle pp \vert \vert rr \rangle)_{4}$'
# Row 8: screen2_1
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 'screen2_1'
$ws.Cells.Item(8, 3).Value = 0.0003517492990990623
$ws.Cells.Item(8, 4).Value = '$(\langle qq \vert \vert ss \rangle)_{1}$'
# Row 9: screen2_2
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 'screen2_2'
$ws.Cells.Item(9, 3).Value = 0.0001793200531106027
$ws.Cells.Item(9, 4).Value = '$(\langle qq \vert \vert ss \rangle)_{2}$'
# Row 10: screen2_3
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 'screen2_3'
$ws.Cells.Item(10, 3).Value = 0.0001475139073027967
$ws.Cells.Item(10, 4).Value = '$(\langle qq \vert \vert ss \rangle)_{3}$'
# Row 11: screen2_4
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 'screen2_4'
$ws.Cells.Item(11, 3).Value = 0.0001841803903263301
$ws.Cells.Item(11, 4).Value = '$(\langle qq \vert \vert ss \rangle)_{4}$'
# Row 12: eijab_1
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 'eijab_1'
$ws.Cells.Item(12, 3).Value = 0.00009138402687976685
$ws.Cells.Item(12, 4).Value = '$(e^{rs}_{pq})_{1}$'
# Row 13: eijab_2
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = 'eijab_2'
$ws.Cells.Item(13, 3).Value = 0.00003218592074617974
$ws.Cells.Item(13, 4).Value = '$(e^{rs}_{pq})_{2}$'
# Row 14: eijab_3
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = 'eijab_3'
$ws.Cells.Item(14, 3).Value = 0.00004403607992682924
$ws.Cells.Item(14, 4).Value = '$(e^{rs}_{pq})_{3}$'
# Row 15: eijab_4
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = 'eijab_4'
$ws.Cells.Item(15, 3).Value = 0.0000179963552426554
$ws.Cells.Item(15, 4).Value = '$(e^{rs}_{pq})_{4}$'
# Row 16: screenvirt_1
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = 'screenvirt_1'
$ws.Cells.Item(16, 3).Value = 0.0001650525225312486
$ws.Cells.Item(16, 4).Value = '$(\langle ss \vert \vert rr \rangle)_{1}$'
# Row 17: screenvirt_2
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = 'screenvirt_2'
$ws.Cells.Item(17, 3).Value = 0.0004471959729984548
$ws.Cells.Item(17, 4).Value = '$(\langle ss \vert \vert rr \rangle)_{2}$'
# Row 18: screenvirt_3
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = 'screenvirt_3'
$ws.Cells.Item(18, 3).Value = 0.0001033953989581824
$ws.Cells.Item(18, 4).Value = '$(\langle ss \vert \vert rr \rangle)_{3}$'
# Row 19: screenvirt_4
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = 'screenvirt_4'
$ws.Cells.Item(19, 3).Value = 0.00003203442550937291
$ws.Cells.Item(19, 4).Value = '$(\langle ss \vert \vert rr \rangle)_{4}$'
# Row 20: Fr1
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = 'Fr1'
$ws.Cells.Item(20, 3).Value = 0.0004010985483623878
$ws.Cells.Item(20, 4).Value = '$(F_{r})_{1}$'
# Row 21: Fr2
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = 'Fr2'
$ws.Cells.Item(21, 3).Value = 0.00002620094815180459
$ws.Cells.Item(21, 4).Value = '$(F_{r})_{2}$'
# Row 22: Fr3
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = 'Fr3'
$ws.Cells.Item(22, 3).Value = 0.0002409746391146895
$ws.Cells.Item(22, 4).Value = '$(F_{r})_{3}$'
# Row 23: Fr4
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = 'Fr4'
$ws.Cells.Item(23, 3).Value = 0.002364669119608438
$ws.Cells.Item(23, 4).Value = '$(F_{r})_{4}$'
# Row 24: Fs1
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = 'Fs1'
$ws.Cells.Item(24, 3).Value = 0.00008886283890297048
$ws.Cells.Item(24, 4).Value = '$(F_{s})_{1}$'
# Row 25: Fs2
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = 'Fs2'
$ws.Cells.Item(25, 3).Value = 0.00001667501288356461
$ws.Cells.Item(25, 4).Value = '$(F_{s})_{2}$'
# Row 26: Fs3
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = 'Fs3'
$ws.Cells.Item(26, 3).Value = 0.00007942505582164792
$ws.Cells.Item(26, 4).Value = '$(F_{s})_{3}$'
# Row 27: Fs4
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = 'Fs4'
$ws.Cells.Item(27, 3).Value = 0.00007366839851843096
$ws.Cells.Item(27, 4).Value = '$(F_{s})_{4}$'
# Row 28: occr1
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = 'occr1'
$ws.Cells.Item(28, 3).Value = 0.00004861308295621949
$ws.Cells.Item(28, 4).Value = '$(\eta_{r})_{1}$'
# Row 29: occr2
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = 'occr2'
$ws.Cells.Item(29, 3).Value = 0.00004840032303823423
$ws.Cells.Item(29, 4).Value = '$(\eta_{r})_{2}$'
# Row 30: occr3
$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = 'occr3'
$ws.Cells.Item(30, 3).Value = 0.0004999630720790536
$ws.Cells.Item(30, 4).Value = '$(\eta_{r})_{3}$'
# Row 31: occr4
$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = 'occr4'
$ws.Cells.Item(31, 3).Value = 0.00006860151238972551
$ws.Cells.Item(31, 4).Value = '$(\eta_{r})_{4}$'
# Row 32: occs1
$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = 'occs1'
$ws.Cells.Item(32, 3).Value = 0.00001137741631172372
$ws.Cells.Item(32, 4).Value = '$(\eta_{s})_{1}$'
# Row 33: occs4
$ws.Cells.Item(33, 1).Value = 34
$ws.Cells.Item(33, 2).Value = 'occs4'
$ws.Cells.Item(33, 3).Value = 0.00001779033012562821
$ws.Cells.Item(33, 4).Value = '$(\eta_{s})_{4}$'
# Row 34: SCFFr1
$ws.Cells.Item(34, 1).Value = 35
$ws.Cells.Item(34, 2).Value = 'SCFFr1'
$ws.Cells.Item(34, 3).Value = 0.0001263330810008015
$ws.Cells.Item(34, 4).Value = '$(F_{r}^{\text{SCF}})_{1}$'
# Row 35: SCFFr2
$ws.Cells.Item(35, 1).Value = 36
$ws.Cells.Item(35, 2).Value = 'SCFFr2'
$ws.Cells.Item(35, 3).Value = 0.00005618717644894813
$ws.Cells.Item(35, 4).Value = '$(F_{r}^{\text{SCF}})_{2}$'
# Row 36: SCFFr3
$ws.Cells.Item(36, 1).Value = 37
$ws.Cells.Item(36, 2).Value = 'SCFFr3'
$ws.Cells.Item(36, 3).Value = 0.006204654913638819
$ws.Cells.Item(36, 4).Value = '$(F_{r}^{\text{SCF}})_{3}$'
# Row 37: SCFFr4
$ws.Cells.Item(37, 1).Value = 38
$ws.Cells.Item(37, 2).Value = 'SCFFr4'
$ws.Cells.Item(37, 3).Value = 0.00006848031965627506
$ws.Cells.Item(37, 4).Value = '$(F_{r}^{\text{SCF}})_{4}$'
# Row 38: SCFFs1
$ws.Cells.Item(38, 1).Value = 39
$ws.Cells.Item(38, 2).Value = 'SCFFs1'
$ws.Cells.Item(38, 3).Value = 0.0001885269457023166
$ws.Cells.Item(38, 4).Value = '$(F_{s}^{\text{SCF}})_{1}$'
# Row 39: SCFFs2
$ws.Cells.Item(39, 1).Value = 40
$ws.Cells.Item(39, 2).Value = 'SCFFs2'
$ws.Cells.Item(39, 3).Value = 0.0000322585281992085
$ws.Cells.Item(39, 4).Value = '$(F_{s}^{\text{SCF}})_{2}$'
# Row 40: SCFFs3
$ws.Cells.Item(40, 1).Value = 41
$ws.Cells.Item(40, 2).Value = 'SCFFs3'
$ws.Cells.Item(40, 3).Value = 0.00007343425400549137
$ws.Cells.Item(40, 4).Value = '$(F_{s}^{\text{SCF}})_{3}$'
# Row 41: SCFFs4
$ws.Cells.Item(41, 1).Value = 42
$ws.Cells.Item(41, 2).Value = 'SCFFs4'
$ws.Cells.Item(41, 3).Value = 0.0001088417929803504
$ws.Cells.Item(41, 4).Value = '$(F_{s}^{\text{SCF}})_{4}$'
# Row 42: hrr1
$ws.Cells.Item(42, 1).Value = 51
$ws.Cells.Item(42, 2).Value = 'hrr1'
$ws.Cells.Item(42, 3).Value = 0.00005047183029967487
$ws.Cells.Item(42, 4).Value = '$(h_{rr})_{1}$'
# Row 43: hrr2
$ws.Cells.Item(43, 1).Value = 52
$ws.Cells.Item(43, 2).Value = 'hrr2'
$ws.Cells.Item(43, 3).Value = 0.0001601035790009357
$ws.Cells.Item(43, 4).Value = '$(h_{rr})_{2}$'
# Row 44: hrr3
$ws.Cells.Item(44, 1).Value = 53
$ws.Cells.Item(44, 2).Value = 'hrr3'
$ws.Cells.Item(44, 3).Value = 0.00008761297843121861
$ws.Cells.Item(44, 4).Value = '$(h_{rr})_{3}$'
# Row 45: hrr4
$ws.Cells.Item(45, 1).Value = 54
$ws.Cells.Item(45, 2).Value = 'hrr4'
$ws.Cells.Item(45, 3).Value = 0.0000980941093322254
$ws.Cells.Item(45, 4).Value = '$(h_{rr})_{4}$'
# Row 46: hss1
$ws.Cells.Item(46, 1).Value = 55
$ws.Cells.Item(46, 2).Value = 'hss1'
$ws.Cells.Item(46, 3).Value = 0.00008448603622286052
$ws.Cells.Item(46, 4).Value = '$(h_{ss})_{1}$'
# Row 47: hss2
$ws.Cells.Item(47, 1).Value = 56
$ws.Cells.Item(47, 2).Value = 'hss2'
$ws.Cells.Item(47, 3).Value = 0.00003421460229764721
$ws.Cells.Item(47, 4).Value = '$(h_{ss})_{2}$'
# Row 48: hss3
$ws.Cells.Item(48, 1).Value = 57
$ws.Cells.Item(48, 2).Value = 'hss3'
$ws.Cells.Item(48, 3).Value = 0.00004637775434270387
$ws.Cells.Item(48, 4).Value = '$(h_{ss})_{3}$'
# Row 49: hss4
$ws.Cells.Item(49, 1).Value = 58
$ws.Cells.Item(49, 2).Value = 'hss4'
$ws.Cells.Item(49, 3).Value = 0.00004152585970472261
$ws.Cells.Item(49, 4).Value = '$(h_{ss})_{4}$'
# Row 50: hpp
$ws.Cells.Item(50, 1).Value = 59
$ws.Cells.Item(50, 2).Value = 'hpp'
$ws.Cells.Item(50, 3).Value = 0.002750638540069508
$ws.Cells.Item(50, 4).Value = '$h_{pp}$'
# Row 51: Fp
$ws.Cells.Item(51, 1).Value = 61
$ws.Cells.Item(51, 2).Value = 'Fp'
$ws.Cells.Item(51, 3).Value = 0.00151652707486342
$ws.Cells.Item(51, 4).Value = '$F_{p}$'
# Row 52: Fq
$ws.Cells.Item(52, 1).Value = 62
$ws.Cells.Item(52, 2).Value = 'Fq'
$ws.Cells.Item(52, 3).Value = 0.0004916753225051912
$ws.Cells.Item(52, 4).Value = '$F_{q}$'
# Row 53: occp
$ws.Cells.Item(53, 1).Value = 63
$ws.Cells.Item(53, 2).Value = 'occp'
$ws.Cells.Item(53, 3).Value = 0.00001835565339829359
$ws.Cells.Item(53, 4).Value = '$\eta_{p}$'
# Row 54: occq
$ws.Cells.Item(54, 1).Value = 64
$ws.Cells.Item(54, 2).Value = 'occq'
$ws.Cells.Item(54, 3).Value = 0.002727938032168187
$ws.Cells.Item(54, 4).Value = '$\eta_{q}$'
# Row 55: SCFFp
$ws.Cells.Item(55, 1).Value = 65
$ws.Cells.Item(55, 2).Value = 'SCFFp'
$ws.Cells.Item(55, 3).Value = 0.0003144343292138878
$ws.Cells.Item(55, 4).Value = '$F_{p}^{\text{SCF}}$'

# Ensure column A data cells (including newly added rows 54-55) carry the same style as the rest of column A
$ws.Range("A2").Copy()
$ws.Range("A2:A55").PasteSpecial(-4122)

$excel.CutCopyMode = 0
